$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update cell values for rows 20-23 (columns A-D) ---
# Row 20
$ws.Range("A20").Value = "ENW000010"
$ws.Range("B20").Value = " OPQA-1968||OPQA-1969 ||OPQA-1970||OPQA-1979||OPQA-1986||OPQA-3864 "
$ws.Range("C20").Value = "Fb login with account setting page linking modal"
$ws.Range("D20").Value = "Y"

# Row 21
$ws.Range("A21").Value = "ENW000012"
$ws.Range("B21").Value = " OPQA-1968||OPQA-1969 ||OPQA-1970||OPQA-1979||OPQA-1986||OPQA-3864 "
$ws.Range("C21").Value = "LI login with account setting page linking modal"
$ws.Range("D21").Value = "Y"

# Row 22
$ws.Range("A22").Value = "ENW00029"
$ws.Range("B22").Value = "OPQA-1919||OPQA-1915"
$ws.Range("C22").Value = "steam login"
$ws.Range("D22").Value = "Y"

# Row 23
$ws.Range("A23").Value = "ENW000011"
$ws.Range("B23").Value = "OPQA-3196"
$ws.Range("C23").Value = "As a user, I want to be able to see all emails that are associated to my Neon identity under the account page"
$ws.Range("D23").Value = "Y"

# --- 2. Normalize formatting to match the rest of the table (thin box border, no special fill) ---
# "No-wrap, bordered" look (same as column A / D throughout the table)
$ws.Range("D16").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D20:D23").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Wrap, bordered" look (same as column B elsewhere in the table)
$ws.Range("B16").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Row heights ---
$ws.Rows.Item(20).RowHeight = 75
$ws.Rows.Item(21).RowHeight = 60

# --- 4. Column B width ---
$ws.Columns.Item(2).ColumnWidth = 24.7

# --- 5. View / selection state ---
$ws.Range("B22").Select()
